$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2
$ws.Range("A2").Value = 5151
$ws.Range("B2").Value = "ريد بل - 250 مل"
$ws.Range("D2").Value = 1065
$ws.Range("E2").Value = "YES"

# Add row 3
$ws.Range("A3").Value = 5152
$ws.Range("B3").Value = "ريد بل فرى شوجر - 250 مل"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 1065
$ws.Range("E3").Value = "YES"

# Add row 4
$ws.Range("A4").Value = 7630
$ws.Range("B4").Value = "فيورى جولد - 400 مل"
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 205
$ws.Range("E4").Value = "YES"
